$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.123091816902161
$ws.Range("B1").Value = 2.270363807678223
$ws.Range("C1").Value = 10.18928050994873
$ws.Range("D1").Value = 1.893985867500305
$ws.Range("E1").Value = 1.287002325057983
